# Apply edits described by the diff:
# - Sheet1 ("Sheet1") gets 4 new rows of company names added, and a
#   pre-existing typo fixed ("Dabus India" -> "Dabur India"), keeping the
#   column alphabetically sorted.
# - The "AllCompanies" sheet content is unchanged (its shared-string
#   indices merely shift because of the newly inserted strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Coal India"
$ws.Range("A2").Value = "Dabur India"
$ws.Range("A3").Value = "Dalmia Bhara"
$ws.Range("A4").Value = "Infosys"
$ws.Range("A5").Value = "KIOCL"
$ws.Range("A6").Value = "Knitworth Export"
$ws.Range("A7").Value = "Maruti Suzuki"
$ws.Range("A8").Value = "Rites"

$ws.Range("A1").Select()
